# Auto-generated from diff plan. Applies every cell-level change
# described by the commit diff across sheets ARM, CUL, LTW, WVR.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 3289.7
$ws.Cells.Item(61, 9).Value = 1374.6666
$ws.Cells.Item(61, 10).Value = 4856.5454
$ws.Cells.Item(61, 11).Value = 1374.6666
$ws.Cells.Item(61, 12).Value = 4856.5454
$ws.Cells.Item(61, 13).Value = -1162.6666
$ws.Cells.Item(61, 14).Value = -5280.5454
$ws.Cells.Item(136, 8).Value = 3289.7
$ws.Cells.Item(136, 9).Value = 1374.6666
$ws.Cells.Item(136, 10).Value = 4856.5454
$ws.Cells.Item(136, 11).Value = 4123.9998
$ws.Cells.Item(136, 12).Value = 14569.6362
$ws.Cells.Item(136, 13).Value = -1573.9998
$ws.Cells.Item(136, 14).Value = -19669.6362
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 99.69
$ws.Cells.Item(2, 9).Value = 0
$ws.Cells.Item(2, 10).Value = 99.69
$ws.Cells.Item(2, 11).Value = 0
$ws.Cells.Item(2, 12).Value = 598.14
$ws.Cells.Item(2, 14).Value = -824.14
$ws.Cells.Item(2, 13).ClearContents()
$ws.Cells.Item(6, 8).Value = 500
$ws.Cells.Item(6, 9).Value = 500
$ws.Cells.Item(6, 10).Value = 0
$ws.Cells.Item(6, 11).Value = 1500
$ws.Cells.Item(6, 12).Value = 0
$ws.Cells.Item(6, 13).Value = -1387
$ws.Cells.Item(40, 8).Value = 1200
$ws.Cells.Item(40, 9).Value = 0
$ws.Cells.Item(40, 10).Value = 1200
$ws.Cells.Item(40, 11).Value = 0
$ws.Cells.Item(40, 12).Value = 4800
$ws.Cells.Item(40, 14).Value = -4938
$ws.Cells.Item(46, 8).Value = 3301
$ws.Cells.Item(46, 9).Value = 2903
$ws.Cells.Item(46, 10).Value = 3500
$ws.Cells.Item(46, 11).Value = 8709
$ws.Cells.Item(46, 12).Value = 10500
$ws.Cells.Item(46, 13).Value = -8618
$ws.Cells.Item(46, 14).Value = -10682
$ws.Cells.Item(68, 8).Value = 3712.5715
$ws.Cells.Item(68, 9).Value = 2000
$ws.Cells.Item(68, 10).Value = 3998
$ws.Cells.Item(68, 11).Value = 6000
$ws.Cells.Item(68, 12).Value = 11994
$ws.Cells.Item(68, 13).Value = -5189
$ws.Cells.Item(68, 14).Value = -13616
$ws.Cells.Item(71, 8).Value = 3712.5715
$ws.Cells.Item(71, 9).Value = 2000
$ws.Cells.Item(71, 10).Value = 3998
$ws.Cells.Item(71, 11).Value = 18000
$ws.Cells.Item(71, 12).Value = 35982
$ws.Cells.Item(71, 13).Value = -13944
$ws.Cells.Item(71, 14).Value = -44094
$ws.Cells.Item(114, 8).Value = 2014
$ws.Cells.Item(114, 9).Value = 2014
$ws.Cells.Item(114, 10).Value = 0
$ws.Cells.Item(114, 11).Value = 6042
$ws.Cells.Item(114, 12).Value = 0
$ws.Cells.Item(114, 13).Value = -2788
$ws.Cells.Item(114, 14).ClearContents()
$ws.Cells.Item(120, 8).ClearContents()
$ws.Cells.Item(120, 9).ClearContents()
$ws.Cells.Item(120, 10).ClearContents()
$ws.Cells.Item(120, 11).ClearContents()
$ws.Cells.Item(120, 12).ClearContents()
$ws.Cells.Item(120, 13).ClearContents()
$ws.Cells.Item(120, 14).ClearContents()
$ws.Cells.Item(121, 8).ClearContents()
$ws.Cells.Item(121, 9).ClearContents()
$ws.Cells.Item(121, 10).ClearContents()
$ws.Cells.Item(121, 11).ClearContents()
$ws.Cells.Item(121, 12).ClearContents()
$ws.Cells.Item(121, 13).ClearContents()
$ws.Cells.Item(121, 14).ClearContents()
$ws.Cells.Item(122, 8).ClearContents()
$ws.Cells.Item(122, 9).ClearContents()
$ws.Cells.Item(122, 10).ClearContents()
$ws.Cells.Item(122, 11).ClearContents()
$ws.Cells.Item(122, 12).ClearContents()
$ws.Cells.Item(122, 13).ClearContents()
$ws.Cells.Item(122, 14).ClearContents()
$ws.Cells.Item(123, 8).ClearContents()
$ws.Cells.Item(123, 9).ClearContents()
$ws.Cells.Item(123, 10).ClearContents()
$ws.Cells.Item(123, 11).ClearContents()
$ws.Cells.Item(123, 12).ClearContents()
$ws.Cells.Item(123, 13).ClearContents()
$ws.Cells.Item(123, 14).ClearContents()
$ws.Cells.Item(124, 8).ClearContents()
$ws.Cells.Item(124, 9).ClearContents()
$ws.Cells.Item(124, 10).ClearContents()
$ws.Cells.Item(124, 11).ClearContents()
$ws.Cells.Item(124, 12).ClearContents()
$ws.Cells.Item(124, 13).ClearContents()
$ws.Cells.Item(124, 14).ClearContents()
$ws.Cells.Item(125, 8).ClearContents()
$ws.Cells.Item(125, 9).ClearContents()
$ws.Cells.Item(125, 10).ClearContents()
$ws.Cells.Item(125, 11).ClearContents()
$ws.Cells.Item(125, 12).ClearContents()
$ws.Cells.Item(125, 13).ClearContents()
$ws.Cells.Item(125, 14).ClearContents()
$ws.Cells.Item(126, 8).ClearContents()
$ws.Cells.Item(126, 9).ClearContents()
$ws.Cells.Item(126, 10).ClearContents()
$ws.Cells.Item(126, 11).ClearContents()
$ws.Cells.Item(126, 12).ClearContents()
$ws.Cells.Item(126, 13).ClearContents()
$ws.Cells.Item(126, 14).ClearContents()
$ws.Cells.Item(127, 8).ClearContents()
$ws.Cells.Item(127, 9).ClearContents()
$ws.Cells.Item(127, 10).ClearContents()
$ws.Cells.Item(127, 11).ClearContents()
$ws.Cells.Item(127, 12).ClearContents()
$ws.Cells.Item(127, 13).ClearContents()
$ws.Cells.Item(127, 14).ClearContents()
$ws.Cells.Item(128, 8).ClearContents()
$ws.Cells.Item(128, 9).ClearContents()
$ws.Cells.Item(128, 10).ClearContents()
$ws.Cells.Item(128, 11).ClearContents()
$ws.Cells.Item(128, 12).ClearContents()
$ws.Cells.Item(128, 13).ClearContents()
$ws.Cells.Item(128, 14).ClearContents()
$ws.Cells.Item(129, 8).ClearContents()
$ws.Cells.Item(129, 9).ClearContents()
$ws.Cells.Item(129, 10).ClearContents()
$ws.Cells.Item(129, 11).ClearContents()
$ws.Cells.Item(129, 12).ClearContents()
$ws.Cells.Item(129, 13).ClearContents()
$ws.Cells.Item(129, 14).ClearContents()
$ws.Cells.Item(130, 8).ClearContents()
$ws.Cells.Item(130, 9).ClearContents()
$ws.Cells.Item(130, 10).ClearContents()
$ws.Cells.Item(130, 11).ClearContents()
$ws.Cells.Item(130, 12).ClearContents()
$ws.Cells.Item(130, 13).ClearContents()
$ws.Cells.Item(130, 14).ClearContents()
$ws.Cells.Item(131, 8).ClearContents()
$ws.Cells.Item(131, 9).ClearContents()
$ws.Cells.Item(131, 10).ClearContents()
$ws.Cells.Item(131, 11).ClearContents()
$ws.Cells.Item(131, 12).ClearContents()
$ws.Cells.Item(131, 13).ClearContents()
$ws.Cells.Item(131, 14).ClearContents()
$ws.Cells.Item(132, 8).ClearContents()
$ws.Cells.Item(132, 9).ClearContents()
$ws.Cells.Item(132, 10).ClearContents()
$ws.Cells.Item(132, 11).ClearContents()
$ws.Cells.Item(132, 12).ClearContents()
$ws.Cells.Item(132, 13).ClearContents()
$ws.Cells.Item(132, 14).ClearContents()
$ws.Cells.Item(133, 8).ClearContents()
$ws.Cells.Item(133, 9).ClearContents()
$ws.Cells.Item(133, 10).ClearContents()
$ws.Cells.Item(133, 11).ClearContents()
$ws.Cells.Item(133, 12).ClearContents()
$ws.Cells.Item(133, 13).ClearContents()
$ws.Cells.Item(133, 14).ClearContents()
$ws.Cells.Item(134, 8).ClearContents()
$ws.Cells.Item(134, 9).ClearContents()
$ws.Cells.Item(134, 10).ClearContents()
$ws.Cells.Item(134, 11).ClearContents()
$ws.Cells.Item(134, 12).ClearContents()
$ws.Cells.Item(134, 13).ClearContents()
$ws.Cells.Item(134, 14).ClearContents()
$ws.Cells.Item(136, 8).ClearContents()
$ws.Cells.Item(136, 9).ClearContents()
$ws.Cells.Item(136, 10).ClearContents()
$ws.Cells.Item(136, 11).ClearContents()
$ws.Cells.Item(136, 12).ClearContents()
$ws.Cells.Item(136, 13).ClearContents()
$ws.Cells.Item(136, 14).ClearContents()
$ws.Cells.Item(137, 8).ClearContents()
$ws.Cells.Item(137, 9).ClearContents()
$ws.Cells.Item(137, 10).ClearContents()
$ws.Cells.Item(137, 11).ClearContents()
$ws.Cells.Item(137, 12).ClearContents()
$ws.Cells.Item(137, 13).ClearContents()
$ws.Cells.Item(137, 14).ClearContents()
$ws.Cells.Item(138, 8).ClearContents()
$ws.Cells.Item(138, 9).ClearContents()
$ws.Cells.Item(138, 10).ClearContents()
$ws.Cells.Item(138, 11).ClearContents()
$ws.Cells.Item(138, 12).ClearContents()
$ws.Cells.Item(138, 13).ClearContents()
$ws.Cells.Item(138, 14).ClearContents()
$ws.Cells.Item(139, 8).ClearContents()
$ws.Cells.Item(139, 9).ClearContents()
$ws.Cells.Item(139, 10).ClearContents()
$ws.Cells.Item(139, 11).ClearContents()
$ws.Cells.Item(139, 12).ClearContents()
$ws.Cells.Item(139, 13).ClearContents()
$ws.Cells.Item(139, 14).ClearContents()
$ws.Cells.Item(140, 8).ClearContents()
$ws.Cells.Item(140, 9).ClearContents()
$ws.Cells.Item(140, 10).ClearContents()
$ws.Cells.Item(140, 11).ClearContents()
$ws.Cells.Item(140, 12).ClearContents()
$ws.Cells.Item(140, 13).ClearContents()
$ws.Cells.Item(140, 14).ClearContents()
$ws.Cells.Item(141, 8).ClearContents()
$ws.Cells.Item(141, 9).ClearContents()
$ws.Cells.Item(141, 10).ClearContents()
$ws.Cells.Item(141, 11).ClearContents()
$ws.Cells.Item(141, 12).ClearContents()
$ws.Cells.Item(141, 13).ClearContents()
$ws.Cells.Item(141, 14).ClearContents()
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(124, 8).Value = 42900
$ws.Cells.Item(124, 9).Value = 0
$ws.Cells.Item(124, 10).Value = 42900
$ws.Cells.Item(124, 11).Value = 0
$ws.Cells.Item(124, 12).Value = 42900
$ws.Cells.Item(124, 14).Value = -52720
$ws.Cells.Item(125, 8).Value = 0
$ws.Cells.Item(125, 9).Value = 0
$ws.Cells.Item(125, 10).Value = 0
$ws.Cells.Item(125, 11).Value = 0
$ws.Cells.Item(125, 12).Value = 0
$ws.Cells.Item(127, 8).Value = 0
$ws.Cells.Item(127, 9).Value = 0
$ws.Cells.Item(127, 10).Value = 0
$ws.Cells.Item(127, 11).Value = 0
$ws.Cells.Item(127, 12).Value = 0
$ws.Cells.Item(128, 8).Value = 0
$ws.Cells.Item(128, 9).Value = 0
$ws.Cells.Item(128, 10).Value = 0
$ws.Cells.Item(128, 11).Value = 0
$ws.Cells.Item(128, 12).Value = 0
$ws.Cells.Item(129, 8).Value = 0
$ws.Cells.Item(129, 9).Value = 0
$ws.Cells.Item(129, 10).Value = 0
$ws.Cells.Item(129, 11).Value = 0
$ws.Cells.Item(129, 12).Value = 0
$ws.Cells.Item(130, 8).Value = 0
$ws.Cells.Item(130, 9).Value = 0
$ws.Cells.Item(130, 10).Value = 0
$ws.Cells.Item(130, 11).Value = 0
$ws.Cells.Item(130, 12).Value = 0
$ws.Cells.Item(131, 8).Value = 0
$ws.Cells.Item(131, 9).Value = 0
$ws.Cells.Item(131, 10).Value = 0
$ws.Cells.Item(131, 11).Value = 0
$ws.Cells.Item(131, 12).Value = 0
$ws.Cells.Item(132, 8).Value = 3922.5
$ws.Cells.Item(132, 9).Value = 3034.125
$ws.Cells.Item(132, 10).Value = 5699.25
$ws.Cells.Item(132, 11).Value = 9102.375
$ws.Cells.Item(132, 12).Value = 17097.75
$ws.Cells.Item(132, 13).Value = -6572.375
$ws.Cells.Item(132, 14).Value = -22157.75
$ws.Cells.Item(133, 8).Value = 0
$ws.Cells.Item(133, 9).Value = 0
$ws.Cells.Item(133, 10).Value = 0
$ws.Cells.Item(133, 11).Value = 0
$ws.Cells.Item(133, 12).Value = 0
$ws.Cells.Item(134, 8).Value = 0
$ws.Cells.Item(134, 9).Value = 0
$ws.Cells.Item(134, 10).Value = 0
$ws.Cells.Item(134, 11).Value = 0
$ws.Cells.Item(134, 12).Value = 0
$ws.Cells.Item(135, 8).Value = 0
$ws.Cells.Item(135, 9).Value = 0
$ws.Cells.Item(135, 10).Value = 0
$ws.Cells.Item(135, 11).Value = 0
$ws.Cells.Item(135, 12).Value = 0
$ws.Cells.Item(136, 8).Value = 1660
$ws.Cells.Item(136, 9).Value = 1660
$ws.Cells.Item(136, 10).Value = 0
$ws.Cells.Item(136, 11).Value = 4980
$ws.Cells.Item(136, 12).Value = 0
$ws.Cells.Item(136, 13).Value = -2430
$ws.Cells.Item(137, 8).Value = 0
$ws.Cells.Item(137, 9).Value = 0
$ws.Cells.Item(137, 10).Value = 0
$ws.Cells.Item(137, 11).Value = 0
$ws.Cells.Item(137, 12).Value = 0
$ws.Cells.Item(138, 8).Value = 0
$ws.Cells.Item(138, 9).Value = 0
$ws.Cells.Item(138, 10).Value = 0
$ws.Cells.Item(138, 11).Value = 0
$ws.Cells.Item(138, 12).Value = 0
$ws.Cells.Item(139, 8).Value = 0
$ws.Cells.Item(139, 9).Value = 0
$ws.Cells.Item(139, 10).Value = 0
$ws.Cells.Item(139, 11).Value = 0
$ws.Cells.Item(139, 12).Value = 0
$ws.Cells.Item(140, 8).Value = 60429
$ws.Cells.Item(140, 9).Value = 0
$ws.Cells.Item(140, 10).Value = 60429
$ws.Cells.Item(140, 11).Value = 0
$ws.Cells.Item(140, 12).Value = 60429
$ws.Cells.Item(140, 14).Value = -70789
$ws.Cells.Item(141, 8).Value = 0
$ws.Cells.Item(141, 9).Value = 0
$ws.Cells.Item(141, 10).Value = 0
$ws.Cells.Item(141, 11).Value = 0
$ws.Cells.Item(141, 12).Value = 0
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(15, 8).Value = 0
$ws.Cells.Item(15, 10).Value = 0
$ws.Cells.Item(15, 12).Value = 0
$ws.Cells.Item(15, 14).ClearContents()
$ws.Cells.Item(18, 8).Value = 0
$ws.Cells.Item(18, 10).Value = 0
$ws.Cells.Item(18, 12).Value = 0
$ws.Cells.Item(18, 14).ClearContents()
$ws.Cells.Item(69, 8).Value = 0
$ws.Cells.Item(69, 10).Value = 0
$ws.Cells.Item(69, 12).Value = 0
$ws.Cells.Item(69, 14).ClearContents()
$ws.Cells.Item(72, 8).Value = 0
$ws.Cells.Item(72, 10).Value = 0
$ws.Cells.Item(72, 12).Value = 0
$ws.Cells.Item(72, 14).ClearContents()
$ws.Cells.Item(132, 8).Value = 1983.9231
$ws.Cells.Item(132, 9).Value = 784.8571
$ws.Cells.Item(132, 10).Value = 3382.8333
$ws.Cells.Item(132, 11).Value = 2354.5713
$ws.Cells.Item(132, 12).Value = 10148.4999
$ws.Cells.Item(132, 13).Value = 175.4287000000004
$ws.Cells.Item(132, 14).Value = -15208.4999
